$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 96142350
$ws.Range("B2").Value = 90653
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 4364
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "Dropptaggsvamp"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "Hydnellum ferrugineum"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "Muggelid, Dls"
$ws.Range("Q2").Value = 333288.6459826281
$ws.Range("R2").Value = 6498947.551675561
$ws.Range("S2").Value = 5
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2021-09-14"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2021-09-14"
$ws.Range("AC2").Value = $null
$ws.Range("AI2").Value = $null
$ws.Range("AO2").Value = $null
$ws.Range("AW2").NumberFormat = "@"
$ws.Range("AW2").Value = "Anton Larsson"
$ws.Range("AX2").NumberFormat = "@"
$ws.Range("AX2").Value = "Anton Larsson"

# Row 3
$ws.Range("A3").Value = 69173232
$ws.Range("B3").Value = 90655
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 788
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "Gul taggsvamp"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "Hydnellum geogenium"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "(Fr.) Banker"
$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value = "Muggelidsbäcken, Dls"
$ws.Range("Q3").Value = 333304.9626084958
$ws.Range("R3").Value = 6498681.251649193
$ws.Range("S3").Value = 10
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2017-08-27"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2017-08-27"
$ws.Range("AI3").NumberFormat = "@"
$ws.Range("AI3").Value = "Äldre barrdominerad skog utmed bäck med delvis kvillande lopp"
$ws.Range("AW3").NumberFormat = "@"
$ws.Range("AW3").Value = "Henrik Weibull"
$ws.Range("AX3").NumberFormat = "@"
$ws.Range("AX3").Value = "Henrik Weibull"
$ws.Range("AY3").NumberFormat = "@"
$ws.Range("AY3").Value = "Åtgärdsprogram för mossor i Västra Götalands län"

# Row 4
$ws.Range("A4").Value = 16044867
$ws.Range("B4").Value = 90671
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 4368
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "Dofttaggsvamp"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "Hydnellum suaveolens"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "(Scop.:Fr.) P. Karst."
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "10"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "fruktkroppar"
$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value = "Buxåskullen, 675 m NNO Balketorp, Dls"
$ws.Range("Q4").Value = 333364.6442548583
$ws.Range("R4").Value = 6498805.502037385
$ws.Range("S4").Value = 25
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2014-07-01"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2014-07-01"
$ws.Range("AH4").NumberFormat = "@"
$ws.Range("AH4").Value = "Granskog"
$ws.Range("AQ4").NumberFormat = "@"
$ws.Range("AQ4").Value = "Kjell Eriksson"
$ws.Range("AR4").NumberFormat = "@"
$ws.Range("AR4").Value = "F1402"
$ws.Range("AU4").NumberFormat = "@"
$ws.Range("AU4").Value = "Rolf-Göran Carlsson"
$ws.Range("AV4").NumberFormat = "@"
$ws.Range("AV4").Value = "2021"

# Row 5
$ws.Range("A5").Value = 96142344
$ws.Range("B5").Value = 90319
$ws.Range("E5").Value = 4769
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "Svavelriska"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "Lactarius scrobiculatus"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "(Scop.:Fr.) Fr."
$ws.Range("P5").NumberFormat = "@"
$ws.Range("P5").Value = "Muggelidsbäcken, Dls"
$ws.Range("Q5").Value = 333265.6097484134
$ws.Range("R5").Value = 6498726.768297401
$ws.Range("S5").Value = 5
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2021-09-14"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2021-09-14"
$ws.Range("AW5").NumberFormat = "@"
$ws.Range("AW5").Value = "Anton Larsson"
$ws.Range("AX5").NumberFormat = "@"
$ws.Range("AX5").Value = "Anton Larsson"

# Row 7
$ws.Range("A7").Value = 69173229
$ws.Range("B7").Value = 92864
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "VU"
$ws.Range("E7").Value = 815
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "Stor skogsbäckmossa"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "Hygrohypnum subeugyrium"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "(Renauld & Cardot) Broth."
$ws.Range("Q7").Value = 333243.743277251
$ws.Range("R7").Value = 6498692.219576385

# Row 8
$ws.Range("A8").Value = 69173228
$ws.Range("B8").Value = 108194
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "Godkänd baserat på observatörens uppgifter"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 219711
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "Sårläka"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "Sanicula europaea"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "L."

# Row 9
$ws.Range("A9").Value = 96142343
$ws.Range("B9").Value = 90319
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "Ovaliderad"
$ws.Range("E9").Value = 4769
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "Svavelriska"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "Lactarius scrobiculatus"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q9").Value = 333251.1676608387
$ws.Range("R9").Value = 6498743.562633296
$ws.Range("S9").Value = 5
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = "2021-09-14"
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = "2021-09-14"
$ws.Range("AI9").Value = $null
$ws.Range("AW9").NumberFormat = "@"
$ws.Range("AW9").Value = "Anton Larsson"
$ws.Range("AX9").NumberFormat = "@"
$ws.Range("AX9").Value = "Anton Larsson"
$ws.Range("AY9").Value = $null

# Row 10
$ws.Range("A10").Value = 2199338
$ws.Range("B10").Value = 108193
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 219711
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "Sårläka"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "Sanicula europaea"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "L."
$ws.Range("I10").Value = $null
$ws.Range("J10").Value = $null
$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value = "Balketorp, 400 m NO om, Dls"
$ws.Range("Q10").Value = 333273.3609825537
$ws.Range("R10").Value = 6498504.655406407
$ws.Range("S10").Value = 50
$ws.Range("Y10").NumberFormat = "@"
$ws.Range("Y10").Value = "1979-04-22"
$ws.Range("AA10").NumberFormat = "@"
$ws.Range("AA10").Value = "1979-04-22"
$ws.Range("AH10").Value = $null
$ws.Range("AQ10").Value = $null
$ws.Range("AR10").Value = $null
$ws.Range("AU10").Value = $null
$ws.Range("AV10").Value = $null

# Row 11
$ws.Range("A11").Value = 2192816
$ws.Range("B11").Value = 104489
$ws.Range("E11").Value = 219686
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "Vätteros"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "Lathraea squamaria"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "L."
$ws.Range("P11").NumberFormat = "@"
$ws.Range("P11").Value = "Balketorp, 400 m NO om, Dls"
$ws.Range("Q11").Value = 333273.3609825537
$ws.Range("R11").Value = 6498504.655406407
$ws.Range("S11").Value = 50
$ws.Range("Y11").NumberFormat = "@"
$ws.Range("Y11").Value = "1982-05-19"
$ws.Range("AA11").NumberFormat = "@"
$ws.Range("AA11").Value = "1982-05-19"
$ws.Range("AW11").NumberFormat = "@"
$ws.Range("AW11").Value = "Kjell Eriksson"
$ws.Range("AX11").NumberFormat = "@"
$ws.Range("AX11").Value = "Kjell Eriksson"

# Row 12
$ws.Range("A12").Value = 2202431
$ws.Range("B12").Value = 108194
$ws.Range("E12").Value = 219711
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "Sårläka"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "Sanicula europaea"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "L."
$ws.Range("P12").NumberFormat = "@"
$ws.Range("P12").Value = "Balketorp, 600 m N-NNV om, Dls"
$ws.Range("Q12").Value = 332860.5715872086
$ws.Range("R12").Value = 6498729.504914329
$ws.Range("S12").Value = 50
$ws.Range("Y12").NumberFormat = "@"
$ws.Range("Y12").Value = "1983-06-10"
$ws.Range("AA12").NumberFormat = "@"
$ws.Range("AA12").Value = "1983-06-10"
$ws.Range("AW12").NumberFormat = "@"
$ws.Range("AW12").Value = "Kjell Eriksson"
$ws.Range("AX12").NumberFormat = "@"
$ws.Range("AX12").Value = "Kjell Eriksson"

# Row 13
$ws.Range("A13").Value = 258538
$ws.Range("B13").Value = 92804
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 782
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "Skirmossa"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "Hookeria lucens"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "(Hedw.) Sm."
$ws.Range("P13").NumberFormat = "@"
$ws.Range("P13").Value = "Muggelidsbäcken, strax N om Balketorp, Dls"
$ws.Range("Q13").Value = 333081.3534867804
$ws.Range("R13").Value = 6498342.744522936
$ws.Range("S13").Value = 10
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("Y13").Value = "1999-05-15"
$ws.Range("AA13").NumberFormat = "@"
$ws.Range("AA13").Value = "1999-05-15"
$ws.Range("AC13").NumberFormat = "@"
$ws.Range("AC13").Value = "Tämligen rikligt längs ca 150 m av bäcken. Datum osäkert."
$ws.Range("AI13").NumberFormat = "@"
$ws.Range("AI13").Value = "I skuggig blandskog"
$ws.Range("AO13").NumberFormat = "@"
$ws.Range("AO13").Value = "På fuktig jord i bäckkant"
$ws.Range("AW13").NumberFormat = "@"
$ws.Range("AW13").Value = "Leif Appelgren"
$ws.Range("AX13").NumberFormat = "@"
$ws.Range("AX13").Value = "Leif Appelgren"
